$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value2 = 1.27
$ws.Range("F3").Value2 = 2.32
$ws.Range("G3").Value2 = 3.25
$ws.Range("H3").Value2 = 2.68
$ws.Range("I3").Value2 = 3.95
$ws.Range("J3").Value2 = 2.84
$ws.Range("K3").Value2 = 4.4
$ws.Range("P3").Value2 = 1.56
$ws.Range("J4").Value2 = 5.5
$ws.Range("R4").Value2 = 1.59
$ws.Range("S4").Value2 = 1.98
$ws.Range("W4").Value2 = 3.25
$ws.Range("K5").Value2 = 4
$ws.Range("Q5").Value2 = 2.1
$ws.Range("F6").Value2 = 1.93
$ws.Range("G6").Value2 = 2.24
$ws.Range("H6").Value2 = 3.9
$ws.Range("Q7").Value2 = 1.84
$ws.Range("H8").Value2 = 1.57
$ws.Range("P8").Value2 = 1.25
$ws.Range("Q8").Value2 = 1.58
$ws.Range("I9").Value2 = 2.42
$ws.Range("K9").Value2 = 3.65
$ws.Range("P9").Value2 = 1.8
$ws.Range("Q9").Value2 = 2.06
$ws.Range("F10").Value2 = 9.4
$ws.Range("H10").Value2 = 1.39
$ws.Range("Q10").Value2 = 1.69
$ws.Range("P11").Value2 = 2.52
$ws.Range("F13").Value2 = 3.4
$ws.Range("G13").Value2 = 3.9
$ws.Range("H13").Value2 = 1.94
$ws.Range("I13").Value2 = 2.14
$ws.Range("J13").Value2 = 4
$ws.Range("K13").Value2 = 5.3
$ws.Range("P13").Value2 = 2.44
$ws.Range("Q13").Value2 = 1.41
$ws.Range("F14").Value2 = 2.3
$ws.Range("G14").Value2 = 2.62
$ws.Range("H14").Value2 = 2.66
$ws.Range("I14").Value2 = 3.1
$ws.Range("J14").Value2 = 3.9
$ws.Range("K14").Value2 = 5
$ws.Range("P14").Value2 = 2.56
$ws.Range("Q14").Value2 = 1.52
$ws.Range("P15").Value2 = 3.4
$ws.Range("F16").Value2 = 3
$ws.Range("K16").Value2 = 3.5
$ws.Range("G19").Value2 = 5
$ws.Range("H19").Value2 = 1.96
$ws.Range("P20").Value2 = 1.71
$ws.Range("U20").Value2 = 1.8
$ws.Range("Z20").Value2 = 42
$ws.Range("AA20").Value2 = 200
$ws.Range("AC20").Value2 = 8
$ws.Range("AD20").Value2 = 23
$ws.Range("AE20").Value2 = 110
$ws.Range("AI20").Value2 = 100
$ws.Range("AL20").Value2 = 55
$ws.Range("F21").Value2 = 2.54
$ws.Range("G21").Value2 = 2.6
$ws.Range("O21").Value2 = 1.44
$ws.Range("Q21").Value2 = 2.34
$ws.Range("U21").Value2 = 1.96
$ws.Range("Y21").Value2 = 11
$ws.Range("AE21").Value2 = 46
$ws.Range("AH21").Value2 = 20
$ws.Range("F22").Value2 = 1.83
$ws.Range("H22").Value2 = 4.1
$ws.Range("I22").Value2 = 5.4
$ws.Range("J22").Value2 = 2.64
$ws.Range("K22").Value2 = 3.7
$ws.Range("P22").Value2 = 1.48
$ws.Range("Q22").Value2 = 2.78
$ws.Range("F23").Value2 = 14
$ws.Range("G23").Value2 = 17.5
$ws.Range("H23").Value2 = 1.27
$ws.Range("P23").Value2 = 2.04
$ws.Range("Q23").Value2 = 1.81
$ws.Range("F24").Value2 = 1.9
$ws.Range("G24").Value2 = 2.02
$ws.Range("H24").Value2 = 4.6
$ws.Range("J24").Value2 = 3.4
$ws.Range("K24").Value2 = 3.75
$ws.Range("F26").Value2 = 3.2
$ws.Range("G26").Value2 = 3.6
$ws.Range("H26").Value2 = 2.46
$ws.Range("I26").Value2 = 2.72
$ws.Range("J26").Value2 = 3
$ws.Range("K26").Value2 = 3.35
$ws.Range("P27").Value2 = 1.41
$ws.Range("F29").Value2 = 2.8
$ws.Range("G29").Value2 = 4.2
$ws.Range("H29").Value2 = 2.3
$ws.Range("I29").Value2 = 3
$ws.Range("J29").Value2 = 2.82
$ws.Range("K29").Value2 = 4.9
$ws.Range("P29").Value2 = 1.5
$ws.Range("Q29").Value2 = 2.32
$ws.Range("J30").Value2 = 3.2
$ws.Range("P30").Value2 = 1.75
$ws.Range("Q30").Value2 = 1.8
